# Add a new "Semester" column (F) that classifies each report-semester row
# (column A) as "Wintersemester" or "Sommersemester".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, styled like the other header cells in row 1 (bold, thin
# border all around, centered horizontally, top-aligned vertically).
$ws.Range("F1").Value = "Semester"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160
$ws.Range("F1").Borders.LineStyle = 1

# Data rows 2-43: classify by parity, matching the Winter/Sommer alternation
# already present in column A (row 2 = "Wintersemester ...", row 3 =
# "Sommersemester ...", and so on).
for ($row = 2; $row -le 43; $row++) {
    if (($row % 2) -eq 0) {
        $ws.Cells.Item($row, 6).Value = "Wintersemester"
    } else {
        $ws.Cells.Item($row, 6).Value = "Sommersemester"
    }
}
